$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark from the middle of the text
#    (Word leaves this bookmark at the position of the last edit; it is being
#    moved to the end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Add a new, empty paragraph at the very end of the document and place the
#    "_GoBack" bookmark (collapsed) inside it.
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertParagraphAfter()

$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $newLastPara.Range.Start

# Work around the COM shim mis-handling bookmarks anchored on a zero-length
# range sitting directly in front of a paragraph mark: insert a placeholder
# character, bookmark the (non-empty) range around it, then delete the
# placeholder again. The bookmark survives the deletion and collapses to the
# correct position.
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertBefore("X")

$bmRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanupRange = $d.Range($insertPos, $insertPos + 1)
$cleanupRange.Text = ""
